# Add 2022-Q3 data
#
# 1) Insert a new row at the top of the "总计" (summary) sheet's data table
#    for the new 2022-Q3 quarter, pushing the existing rows down by one.
# 2) Insert a brand-new "2022-Q3" worksheet (positioned right after "总计",
#    i.e. the new second tab) holding the per-fund holdings detail for the
#    new quarter.

$wb = $excel.ActiveWorkbook

# Helper: write a value into a Range as literal text, without Excel's
# "looks like a number" auto-coercion (which would turn "005112" into 5112,
# or "3.30" into 3.3) and without leaving the "number stored as text" quote
# -prefix marker (and its accompanying style) that a plain leading apostrophe
# would add. Building the text via a throw-away formula and then freezing it
# with a Paste-Values keeps the cell's style untouched (matches the source
# data, which was written by a non-Excel tool as plain text).
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '=""&"' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert the 2022-Q3 row above the existing data.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Copy row 6's formatting (incl. style "s=3" on column A) down into the new
# row 7 before shifting data, so the newly-exposed row 7 is styled like its
# siblings.
$summary.Range("A6:D6").Copy($summary.Range("A7:D7"))

# Shift the existing B/C/D values down one row (bottom-up so we don't clobber
# data we still need to read). Column A is left alone -- it is just the
# sequential 0-based row index and already has the right value in every row
# except the newly created row 7.
for ($r = 6; $r -ge 2; $r--) {
    $dest = $r + 1
    $summary.Range("B$dest").Value = $summary.Range("B$r").Value()
    $summary.Range("C$dest").Value = $summary.Range("C$r").Value()
    $summary.Range("D$dest").Value = $summary.Range("D$r").Value()
}
$summary.Range("A7").Value = 5

# Write the new 2022-Q3 summary row.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 0.96

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet, positioned right after "总计".
# ---------------------------------------------------------------------
# Duplicate the existing "2022-Q2" sheet -- it already has the right header
# row/column styling -- then rename it and overwrite its data.
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($null, $summary)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The template has 6 data rows (A1:H6); the new quarter needs 7 (A1:H8).
# Copy the last data row's formatting down to extend the table.
$q3.Range("A6:H6").Copy($q3.Range("A7:H7"))
$q3.Range("A6:H6").Copy($q3.Range("A8:H8"))

# Row 2
$q3.Range("A2").Value = 0
Set-TextValue $q3.Range("B2") "501011"
Set-TextValue $q3.Range("C2") "汇添富中证中药指数（LOF）A"
Set-TextValue $q3.Range("D2") "11.36"
Set-TextValue $q3.Range("E2") "94.73"
Set-TextValue $q3.Range("F2") "3.18"
Set-TextValue $q3.Range("G2") "0.3612"
$q3.Range("H2").Value = 8

# Row 3
$q3.Range("A3").Value = 1
Set-TextValue $q3.Range("B3") "501012"
Set-TextValue $q3.Range("C3") "汇添富中证中药指数（LOF）C"
Set-TextValue $q3.Range("D3") "6.42"
Set-TextValue $q3.Range("E3") "94.73"
Set-TextValue $q3.Range("F3") "3.18"
Set-TextValue $q3.Range("G3") "0.2042"
$q3.Range("H3").Value = 8

# Row 4
$q3.Range("A4").Value = 2
Set-TextValue $q3.Range("B4") "159647"
Set-TextValue $q3.Range("C4") "鹏华中证中药ETF"
Set-TextValue $q3.Range("D4") "6.16"
Set-TextValue $q3.Range("E4") "94.79"
Set-TextValue $q3.Range("F4") "3.13"
Set-TextValue $q3.Range("G4") "0.1928"
$q3.Range("H4").Value = 8

# Row 5
$q3.Range("A5").Value = 3
Set-TextValue $q3.Range("B5") "562390"
Set-TextValue $q3.Range("C5") "银华中证中药ETF"
Set-TextValue $q3.Range("D5") "2.34"
Set-TextValue $q3.Range("E5") "98.09"
Set-TextValue $q3.Range("F5") "3.30"
Set-TextValue $q3.Range("G5") "0.0772"
$q3.Range("H5").Value = 8

# Row 6
$q3.Range("A6").Value = 4
Set-TextValue $q3.Range("B6") "561510"
Set-TextValue $q3.Range("C6") "华泰柏瑞中证中药ETF"
Set-TextValue $q3.Range("D6") "2.02"
Set-TextValue $q3.Range("E6") "95.98"
Set-TextValue $q3.Range("F6") "3.22"
Set-TextValue $q3.Range("G6") "0.0650"
$q3.Range("H6").Value = 8

# Row 7
$q3.Range("A7").Value = 5
Set-TextValue $q3.Range("B7") "005112"
Set-TextValue $q3.Range("C7") "银华中证全指医药卫生指数增强"
Set-TextValue $q3.Range("D7") "1.35"
Set-TextValue $q3.Range("E7") "90.85"
Set-TextValue $q3.Range("F7") "4.08"
Set-TextValue $q3.Range("G7") "0.0551"
$q3.Range("H7").Value = 9

# Row 8
$q3.Range("A8").Value = 6
Set-TextValue $q3.Range("B8") "010487"
Set-TextValue $q3.Range("C8") "中银顺盈回报一年持有期混合"
Set-TextValue $q3.Range("D8") "0.82"
Set-TextValue $q3.Range("E8") "25.62"
Set-TextValue $q3.Range("F8") "0.59"
Set-TextValue $q3.Range("G8") "0.0048"
$q3.Range("H8").Value = 8
